$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = "Presented by,`rBharath Muthuswamy Paran."
Write-Output $tr.Text
